# Rename the "Scanner" sheet to "Session"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Session"

# Delete row 2 (data row) so only the header row remains.
# This also shrinks the worksheet's used-range dimension from A1:F2 to A1:F1.
$ws.Rows.Item(2).Delete()

# Re-assert the "numbers stored as text" ignored-error range so it tracks
# the header row only (A1:F1) instead of the old A1:F2 range.
$ws.Range("A1:F1").Errors.Item(9).Ignore = $true
